# Update "想去人数" (interested-people count) figures in F column
# across the "展览" (sheet1), "演出" (sheet2), and "全部类型" (sheet4) sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 815
$ws1.Range("F4").Value = 1139
$ws1.Range("F9").Value = 393
$ws1.Range("F15").Value = 12949
$ws1.Range("F17").Value = 5312
$ws1.Range("F18").Value = 5540

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 147

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 815
$ws4.Range("F4").Value = 1139
$ws4.Range("F9").Value = 393
$ws4.Range("F15").Value = 12949
$ws4.Range("F16").Value = 147
$ws4.Range("F19").Value = 5312
$ws4.Range("F20").Value = 5540
